# Generate Report for Handoff
# Refresh the "latest handoff" timestamps for the file that was just
# re-handed-off (d88e27aa-c46e-4436-a424-97a4e33712ca.md), on the
# per-language sheets as well as the roll-up Overview sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: "Latest Handoff Datetime" (column E) for that file's row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-19 03:19:03"

# --- de-de sheet: "Latest Handoff Datetime" (column E) for that file's row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-19 03:19:11"

# --- Overview sheet: "Latest Handoff Date" (column D) for that file's row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-19 03:19:11"
